$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename H1, insert new Std_* headers I1:N1, move Obs_Prob to O1 ---
$headers = New-Object 'object[,]' 1,15
$headers[0,0] = '#_Agents'
$headers[0,1] = 'Coverage'
$headers[0,2] = 'Avg_Total_Rounds'
$headers[0,3] = 'Avg_Expl_Cost'
$headers[0,4] = 'Avg_Expl_Eff'
$headers[0,5] = 'Avg_Round_Time'
$headers[0,6] = 'Avg_Agent_Step_Time'
$headers[0,7] = 'Avg_Experiment_Time'
$headers[0,8] = 'Std_Total_Rounds'
$headers[0,9] = 'Std_Expl_Cost'
$headers[0,10] = 'Std_Expl_Eff'
$headers[0,11] = 'Std_Round_Time'
$headers[0,12] = 'Std_Agent_Step_Time'
$headers[0,13] = 'Std_Experiment_Time'
$headers[0,14] = 'Obs_Prob'
$ws.Range("A1:O1").Value = $headers

# Apply the bold/bordered header style to the newly-added header cells (K1:O1)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("I1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Data rows 2-13: refreshed experiment averages/std-devs (full recompute) ---
$data = New-Object 'object[,]' 12,15
$data[0,0] = 1.0
$data[0,1] = 1.0
$data[0,2] = 52.764
$data[0,3] = 52.764
$data[0,4] = 3.26659022
$data[0,5] = 0.14041076
$data[0,6] = 0.14041076
$data[0,7] = 7.408108120000001
$data[0,8] = 7.210567840132587
$data[0,9] = 7.210567840132587
$data[0,10] = 0.4702148447944638
$data[0,11] = 0.001014295328010955
$data[0,12] = 0.001014295328010955
$data[0,13] = 1.010627676857371
$data[0,14] = 0.15
$data[1,0] = 1.0
$data[1,1] = 1.0
$data[1,2] = 87.964
$data[1,3] = 87.964
$data[1,4] = 1.9634945
$data[1,5] = 0.1314479
$data[1,6] = 0.1314479
$data[1,7] = 11.54657602
$data[1,8] = 12.88015919940448
$data[1,9] = 12.88015919940448
$data[1,10] = 0.2950596548316437
$data[1,11] = 0.002674962998639079
$data[1,12] = 0.002674962998639079
$data[1,13] = 1.594990560069426
$data[1,14] = 0.85
$data[2,0] = 2.0
$data[2,1] = 1.0
$data[2,2] = 28.972
$data[2,3] = 57.944
$data[2,4] = 2.9974748
$data[2,5] = 0.1681331
$data[2,6] = 0.0840664
$data[2,7] = 2.4351256
$data[2,8] = 4.965130516117795
$data[2,9] = 9.93026103223559
$data[2,10] = 0.4835505048601034
$data[2,11] = 0.002038212847871009
$data[2,12] = 0.00101914569905967
$data[2,13] = 0.4161103568728297
$data[2,14] = 0.15
$data[3,0] = 2.0
$data[3,1] = 1.0
$data[3,2] = 48.418
$data[3,3] = 96.828
$data[3,4] = 1.80582288
$data[3,5] = 0.13557258
$data[3,6] = 0.06778618
$data[3,7] = 3.27226822
$data[3,8] = 9.083628162865452
$data[3,9] = 18.16254861919685
$data[3,10] = 0.3343094677045497
$data[3,11] = 0.003567870140332543
$data[3,12] = 0.001783999584258154
$data[3,13] = 0.5627812265815347
$data[3,14] = 0.85
$data[4,0] = 4.0
$data[4,1] = 1.0
$data[4,2] = 15.094
$data[4,3] = 60.37
$data[4,4] = 2.96269876
$data[4,5] = 0.2033019
$data[4,6] = 0.05082556000000001
$data[4,7] = 0.7669267999999999
$data[4,8] = 3.811853151538998
$data[4,9] = 15.24657391127804
$data[4,10] = 0.680048548493325
$data[4,11] = 0.00432772769985453
$data[4,12] = 0.001081876703142093
$data[4,13] = 0.1934447108588874
$data[4,14] = 0.15
$data[5,0] = 4.0
$data[5,1] = 1.0
$data[5,2] = 25.474
$data[5,3] = 101.856
$data[5,4] = 1.72930402
$data[5,5] = 0.1459651
$data[5,6] = 0.03649126
$data[5,7] = 0.9249828
$data[5,8] = 5.299953113567329
$data[5,9] = 21.18000115433886
$data[5,10] = 0.3515591729518842
$data[5,11] = 0.006068311214166918
$data[5,12] = 0.00151709451384146
$data[5,13] = 0.1717424350139739
$data[5,14] = 0.85
$data[6,0] = 6.0
$data[6,1] = 1.0
$data[6,2] = 10.218
$data[6,3] = 61.308
$data[6,4] = 2.89078352
$data[6,5] = 0.23834036
$data[6,6] = 0.0397233
$data[6,7] = 0.40577338
$data[6,8] = 2.371016541249292
$data[6,9] = 14.22609924749575
$data[6,10] = 0.6120222740740202
$data[6,11] = 0.005137192634769838
$data[6,12] = 0.0008559320939765012
$data[6,13] = 0.09406166487775319
$data[6,14] = 0.15
$data[7,0] = 6.0
$data[7,1] = 1.0
$data[7,2] = 17.564
$data[7,3] = 105.278
$data[7,4] = 1.6924995
$data[7,5] = 0.14907528
$data[7,6] = 0.0248458
$data[7,7] = 0.43350276
$data[7,8] = 4.221557425191369
$data[7,9] = 25.19210193691666
$data[7,10] = 0.3847808865477272
$data[7,11] = 0.008179344807780703
$data[7,12] = 0.001363166168890652
$data[7,13] = 0.09329373236278393
$data[7,14] = 0.85
$data[8,0] = 8.0
$data[8,1] = 1.0
$data[8,2] = 7.564
$data[8,3] = 60.512
$data[8,4] = 2.94196972
$data[8,5] = 0.25918398
$data[8,6] = 0.03239788
$data[8,7] = 0.24491212
$data[8,8] = 1.752725904467776
$data[8,9] = 14.02180723574221
$data[8,10] = 0.6853156447516753
$data[8,11] = 0.006976075143547183
$data[8,12] = 0.0008722030623449609
$data[8,13] = 0.05633992154614036
$data[8,14] = 0.15
$data[9,0] = 8.0
$data[9,1] = 1.0
$data[9,2] = 13.63
$data[9,3] = 108.886
$data[9,4] = 1.65509722
$data[9,5] = 0.14451536
$data[9,6] = 0.01806436
$data[9,7] = 0.2445218
$data[9,8] = 3.669336187393058
$data[9,9] = 29.17216197110074
$data[9,10] = 0.4160260929616332
$data[9,11] = 0.00941751017340211
$data[9,12] = 0.001177250877960734
$data[9,13] = 0.06065823147150353
$data[9,14] = 0.85
$data[10,0] = 10.0
$data[10,1] = 1.0
$data[10,2] = 6.028
$data[10,3] = 60.28
$data[10,4] = 2.99656546
$data[10,5] = 0.2867428599999999
$data[10,6] = 0.02867418
$data[10,7] = 0.17271412
$data[10,8] = 1.587531578591762
$data[10,9] = 15.87531578591762
$data[10,10] = 0.775573309138785
$data[10,11] = 0.00829669519282499
$data[10,12] = 0.0008295846847056643
$data[10,13] = 0.04530979411906243
$data[10,14] = 0.15
$data[11,0] = 10.0
$data[11,1] = 1.0
$data[11,2] = 10.944
$data[11,3] = 109.168
$data[11,4] = 1.66625822
$data[11,5] = 0.1425948
$data[11,6] = 0.0142594
$data[11,7] = 0.15489382
$data[11,8] = 3.137603082961965
$data[11,9] = 31.09974540617105
$data[11,10] = 0.4459830293074974
$data[11,11] = 0.009143890946598532
$data[11,12] = 0.000914416510737979
$data[11,13] = 0.04080276215167317
$data[11,14] = 0.85
$ws.Range("A2:O13").Value = $data
